# Adds the 15 newest Lotofácil draws (rows 403-417) to the bottom of the
# results table, then leaves the selection on G404 (matching the author's
# last-touched cell) as in the committed workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(3473,2,3,4,5,6,7,9,12,13,14,17,18,19,23,25),
    @(3474,1,4,6,8,9,10,12,13,14,15,18,19,20,21,25),
    @(3475,1,3,4,5,6,11,14,15,18,20,21,22,23,24,25),
    @(3476,3,5,6,7,8,10,12,15,16,17,19,20,23,24,25),
    @(3477,2,3,5,6,8,12,15,16,17,18,20,21,23,24,25),
    @(3478,1,4,5,6,7,8,9,10,12,14,15,17,21,22,24),
    @(3479,1,5,6,8,9,10,12,13,15,16,18,19,23,24,25),
    @(3480,3,5,6,8,9,12,13,14,15,16,17,20,21,22,23),
    @(3481,1,2,3,4,5,7,9,10,13,14,19,21,22,23,24),
    @(3482,1,2,4,5,9,10,13,14,15,17,19,20,22,24,25),
    @(3483,2,4,7,8,9,12,13,14,15,17,18,20,22,24,25),
    @(3484,1,3,4,5,7,8,9,11,13,14,15,17,23,24,25),
    @(3485,2,3,4,5,6,9,10,11,12,13,15,18,20,22,25),
    @(3486,2,6,7,8,9,10,11,12,16,17,18,19,22,24,25),
    @(3487,2,4,5,7,10,12,13,15,16,17,18,19,21,22,23)
)

$startRow = 403
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($col = 1; $col -le $values.Length; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}

$ws.Range("G404").Select() | Out-Null
